# Apply the two textual edits described by the diff:
#  1. Slide 11 ("Aside: Singleton Design Pattern"), TextShape 2:
#       "	private " -> "	private static "
#  2. Slide 4 ("Usage"), TextShape 2, the first marioTween.animate(...) line:
#       merge the ", 2.0, " and "500);" runs into a single ", 2.0, 500);" run

$p = $ppt.ActivePresentation

# --- Edit 1: slide 11, singleton pattern "private" field ---
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$tr11 = $sh11.TextFrame.TextRange

# paragraph 6 is "	private MyClass instance;"
$para11 = $tr11.Paragraphs(6, 1)
$run11 = $para11.Runs(1, 1)
$run11.Text = "`tprivate static "

# --- Edit 2: slide 4, Tween.animate(SCALE_X, 2.0, 500); call ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

# paragraph 8 is "marioTween.animate(TweenableParams.SCALE_X, 2.0, 500);"
# runs: 1="marioTween.animate" 2="(" 3="TweenableParams.SCALE_X" 4=", 2.0, " 5="500);"
# Merge runs 4 and 5 into a single run with text ", 2.0, 500);" (keeping run 4's
# formatting, dropping run 5 and its run properties entirely).
$para4 = $tr4.Paragraphs(8, 1)
$run4 = $para4.Runs(4, 1)
$run5 = $para4.Runs(5, 1)
$mergedLength = $run4.Text.Length + $run5.Text.Length
$combined = $tr4.Characters($run4.Start, $mergedLength)
$combined.Text = ", 2.0, 500);"
